# Apply cryptos list update (GitHub Actions scheduled refresh).
# Column D ("Price") values are plain text (e.g. "37.318.69", "1.00") that
# must stay text -- mark each touched D cell as Text ("@") before assigning
# so Excel does not silently coerce them to numbers and drop trailing zeros
# or thousands-separator dots.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "37.318.69"
$ws.Range("E2").Value = "  -1.34%  "
# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.049.55"
$ws.Range("E3").Value = "  -1.51%  "
# Row 4
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$ws.Range("E4").Value = "  -0.10%  "
# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "228.85"
$ws.Range("E5").Value = "  -2.11%  "
# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.614"
$ws.Range("E6").Value = "  -1.75%  "
# Row 7
$ws.Range("E7").Value = "  +0.03%  "
# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "57.18"
$ws.Range("E8").Value = "  -2.40%  "
# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.386"
$ws.Range("E9").Value = "  -1.89%  "
# Row 11
$ws.Range("E11").Value = "  -1.97%  "
# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "14.69"
$ws.Range("E12").Value = "  -2.29%  "
# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "2.351.23"
$ws.Range("E13").Value = "  -1.56%  "
# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "20.80"
$ws.Range("E14").Value = "  -2.61%  "
# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.760"
$ws.Range("E15").Value = "  -2.99%  "
# Row 16
$ws.Range("E16").Value = "  -1.86%  "
# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.046.76"
$ws.Range("E17").Value = "  -1.32%  "
# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "37.222.42"
$ws.Range("E18").Value = "  -1.55%  "
# Row 19
$ws.Range("E19").Value = "  -0.60%  "
# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "69.83"
$ws.Range("E20").Value = "  -2.24%  "
# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0835"
$ws.Range("E21").Value = "  -0.76%  "
# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "226.44"
$ws.Range("E22").Value = "  -1.86%  "
# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$ws.Range("E23").Value = "  -0.08%  "
# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.40"
$ws.Range("E24").Value = "  +0.27%  "
# Row 25
$ws.Range("E25").Value = "  -5.14%  "
# Row 26
$ws.Range("E26").Value = "  -3.41%  "
# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "168.34"
$ws.Range("E27").Value = "  -2.18%  "
# Row 28
$ws.Range("E28").Value = "  -4.74%  "
# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.39"
$ws.Range("E29").Value = "  -1.59%  "
# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "18.95"
$ws.Range("E30").Value = "  -2.85%  "
# Row 31
$ws.Range("E31").Value = "  -2.80%  "
# Row 32
$ws.Range("E32").Value = "  -4.06%  "
# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.62"
$ws.Range("E33").Value = "  -1.82%  "
# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.0615"
$ws.Range("E34").Value = "  -3.11%  "
# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.42"
$ws.Range("E35").Value = "  -2.15%  "
# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.83"
$ws.Range("E36").Value = "  +0.64%  "
# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$ws.Range("E37").Value = "  +0.03%  "
# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "3.25"
$ws.Range("E38").Value = "  -4.55%  "
# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "5.37"
$ws.Range("E39").Value = "  -2.14%  "
# Row 40
$ws.Range("E40").Value = "  -5.22%  "
# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "17.17"
$ws.Range("E41").Value = "  +2.32%  "
# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.493.37"
$ws.Range("E42").Value = "  +2.94%  "
# Row 43
$ws.Range("E43").Value = "  -1.39%  "
# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "96.65"
$ws.Range("E44").Value = "  -5.58%  "
# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0942"
$ws.Range("E45").Value = "  -3.27%  "
# Row 46
$ws.Range("E46").Value = "  +0.71%  "
# Row 47
$ws.Range("E47").Value = "  -4.21%  "
# Row 48
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "7.16"
$ws.Range("E48").Value = "  -2.32%  "
# Row 49
$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "3.87"
$ws.Range("E49").Value = "  -5.77%  "
# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.91"
$ws.Range("E50").Value = "  -2.84%  "
# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.239.39"
$ws.Range("E51").Value = "  -1.49%  "
